# "Generate Report for handback" — mark the two tracked files as handed
# back (in sync with en-US) for both locales, stamp the handback
# datetime, and record the new "Latest Target File" / "Latest Handback
# File" hyperlinks (columns E/F) that the handback report adds.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: B2/C2 (row for ce5c86f6-...md) and B3/C3 (row for
# ffff00af8289-...md) mirror the same status string.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusText
$wsZh.Range("B3").Value = $statusText

$wsZh.Range("G2").Value = "2016-01-26 06:52:59"
$wsZh.Range("G3").Value = "2016-01-26 06:52:59"

$zhMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/f919496bd9d9416b1351c2e3211b030944472595/e2e/ce5c86f6-c426-418d-87ed-595e7ca3b99c.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4994edb14ff70c31437d89cc8501eacb442eb4e9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhMdUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $zhMdUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlfUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f919496bd9d9416b1351c2e3211b030944472595/e2e/ffff00af8289-a5b5-4788-9332-6a0919c53a09.md", "", "", "ffff00af8289-a5b5-4788-9332-6a0919c53a09.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $zhXlfUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $zhMdUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlfUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f919496bd9d9416b1351c2e3211b030944472595/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusText
$wsDe.Range("B3").Value = $statusText

$wsDe.Range("G2").Value = "2016-01-26 06:53:17"
$wsDe.Range("G3").Value = "2016-01-26 06:53:17"

$deMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/f919496bd9d9416b1351c2e3211b030944472595/e2e/ce5c86f6-c426-418d-87ed-595e7ca3b99c.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/87961288df215b1c1253e0908419ff262a2f7c36/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deMdUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $deMdUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlfUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f919496bd9d9416b1351c2e3211b030944472595/e2e/ffff00af8289-a5b5-4788-9332-6a0919c53a09.md", "", "", "ffff00af8289-a5b5-4788-9332-6a0919c53a09.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $deXlfUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $deMdUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlfUrl, "", "", "ce5c86f6-c426-418d-87ed-595e7ca3b99c.9042bc0ba658706f9da0dd94cfdb9860e708325e.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f919496bd9d9416b1351c2e3211b030944472595/.localization-config", "", "", ".localization-config")
